$wb = $excel.ActiveWorkbook

# --- Add the new "Assay" worksheet at the end of the sheet tabs ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$assay = $wb.Worksheets.Add($null, $lastSheet)
$assay.Name = "Assay"

# Populate the "Assay" sheet with the new Settings/Value table
$assay.Range("A1").Value = "Setting"
$assay.Range("B1").Value = "Value"

$assay.Range("A2").Value = "DMSO Tolerance"
$assay.Range("B2").Value = 0.005

$assay.Range("A3").Value = "Well Volume (µL)"
$assay.Range("B3").Value = 25

$assay.Range("A4").Value = "Backfill (µL)"
$assay.Range("B4").Value = 10

$assay.Range("A5").Value = "Allowed Error"
$assay.Range("B5").Value = 0.1

$assay.Range("A6").Value = "Destination Replicates"
$assay.Range("B6").Value = 1

$assay.Range("A7").Value = "Use Intermediate Plates"
$assay.Range("B7").Value = 1

$assay.Range("A8").Value = "DMSO Normalization"
$assay.Range("B8").Value = 1

# Leave the Assay sheet's cursor sitting on H13, matching the authored file
$assay.Range("H13").Select() | Out-Null

# --- Compounds sheet: strip the (duplicate/unused) explicit style off the header row ---
$compounds = $wb.Worksheets.Item("Compounds")
$compounds.Range("A1:F1").ClearFormats() | Out-Null

# --- Patterns becomes the active tab, with E10 selected ---
$patterns = $wb.Worksheets.Item("Patterns")
$patterns.Activate() | Out-Null
$patterns.Range("E10").Select() | Out-Null
